# Scheduled-runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across several Leve rows
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 7000
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H26").Value = 35508
$ws.Range("I26").Value = 1001
$ws.Range("J26").Value = 70015
$ws.Range("K26").Value = 1001
$ws.Range("L26").Value = 70015
$ws.Range("M26").Value = -657
$ws.Range("N26").Value = -70703

$ws.Range("H92").Value = 1573.5834
$ws.Range("I92").Value = 1401.3334
$ws.Range("J92").Value = 1745.8334
$ws.Range("K92").Value = 1401.3334
$ws.Range("L92").Value = 1745.8334
$ws.Range("M92").Value = -153.3334
$ws.Range("N92").Value = -4241.8334

$ws.Range("H131").Value = 2677.5454
$ws.Range("I131").Value = 2004.2354
$ws.Range("J131").Value = 4966.8
$ws.Range("K131").Value = 6012.706200000001
$ws.Range("L131").Value = 14900.4
$ws.Range("M131").Value = -972.7062000000005
$ws.Range("N131").Value = -24980.4

$ws.Range("H138").Value = 3049.423
$ws.Range("I138").Value = 1533.95
$ws.Range("J138").Value = 8101
$ws.Range("K138").Value = 4601.85
$ws.Range("L138").Value = 24303
$ws.Range("M138").Value = 538.1499999999996
$ws.Range("N138").Value = -34583

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 62506770
$ws.Range("I132").Value = 125009350
$ws.Range("J132").Value = 4187.25
$ws.Range("K132").Value = 375028050
$ws.Range("L132").Value = 12561.75
$ws.Range("M132").Value = -375025520
$ws.Range("N132").Value = -17621.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 960459.5600000001
$ws.Range("I86").Value = 1328451.8
$ws.Range("K86").Value = 1328451.8
$ws.Range("M86").Value = -1327328.8

$ws.Range("H89").Value = 960459.5600000001
$ws.Range("I89").Value = 1328451.8
$ws.Range("K89").Value = 6642259
$ws.Range("M89").Value = -6636643

$ws.Range("H132").Value = 39800
$ws.Range("J132").Value = 39800
$ws.Range("L132").Value = 39800
$ws.Range("N132").Value = -49920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 29866.666
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 29866.666
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 29866.666
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -31116.666

$ws.Range("H51").Value = 19366.666
$ws.Range("J51").Value = 19366.666
$ws.Range("L51").Value = 19366.666
$ws.Range("N51").Value = -20838.666

$ws.Range("H59").Value = 24237.5
$ws.Range("I59").Value = 14000
$ws.Range("J59").Value = 27650
$ws.Range("K59").Value = 14000
$ws.Range("L59").Value = 27650
$ws.Range("M59").Value = -12855
$ws.Range("N59").Value = -29940

$ws.Range("H60").Value = 18525
$ws.Range("I60").Value = 16000
$ws.Range("J60").Value = 19366.666
$ws.Range("K60").Value = 16000
$ws.Range("L60").Value = 19366.666
$ws.Range("M60").Value = -15489
$ws.Range("N60").Value = -20388.666

$ws.Range("H61").Value = 19366.666
$ws.Range("J61").Value = 19366.666
$ws.Range("L61").Value = 19366.666
$ws.Range("N61").Value = -20062.666

$ws.Range("H62").Value = 4354.2
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4354.2
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4354.2
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -5602.2

$ws.Range("H65").Value = 4354.2
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4354.2
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 21771
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -28011

$ws.Range("H68").Value = 43200
$ws.Range("J68").Value = 43200
$ws.Range("L68").Value = 43200
$ws.Range("N68").Value = -44698

$ws.Range("H71").Value = 43200
$ws.Range("J71").Value = 43200
$ws.Range("L71").Value = 129600
$ws.Range("N71").Value = -137088

$ws.Range("H74").Value = 19742.285
$ws.Range("J74").Value = 19742.285
$ws.Range("L74").Value = 19742.285
$ws.Range("N74").Value = -21490.285

$ws.Range("H77").Value = 19742.285
$ws.Range("J77").Value = 19742.285
$ws.Range("L77").Value = 59226.855
$ws.Range("N77").Value = -67962.855

$ws.Range("H134").Value = 4732.8887
$ws.Range("I134").Value = 2649.75
$ws.Range("J134").Value = 6399.4
$ws.Range("K134").Value = 7949.25
$ws.Range("L134").Value = 19198.2
$ws.Range("M134").Value = -5414.25
$ws.Range("N134").Value = -24268.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1192.2954
$ws.Range("I107").Value = 721.4
$ws.Range("J107").Value = 1330.7941
$ws.Range("K107").Value = 2164.2
$ws.Range("L107").Value = 3992.3823
$ws.Range("M107").Value = -244.1999999999998
$ws.Range("N107").Value = -7832.3823

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 55004.5
$ws.Range("I18").Value = 10000
$ws.Range("J18").Value = 70006
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 70006
$ws.Range("M18").Value = -9707
$ws.Range("N18").Value = -70592

$ws.Range("H122").Value = 4011.1765
$ws.Range("I122").Value = 2354.6667
$ws.Range("J122").Value = 5874.75
$ws.Range("K122").Value = 7064.000100000001
$ws.Range("L122").Value = 17624.25
$ws.Range("M122").Value = -4614.000100000001
$ws.Range("N122").Value = -22524.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 298
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 298
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 298
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -578

$ws.Range("H30").Value = 1009.75
$ws.Range("I30").Value = 1009.75
$ws.Range("K30").Value = 1009.75
$ws.Range("M30").Value = -901.75

$ws.Range("H122").Value = 2385.0278
$ws.Range("I122").Value = 2095.5652
$ws.Range("J122").Value = 2897.1538
$ws.Range("K122").Value = 6286.6956
$ws.Range("L122").Value = 8691.4614
$ws.Range("M122").Value = -3836.6956
$ws.Range("N122").Value = -13591.4614

$ws.Range("H132").Value = 2720.2122
$ws.Range("I132").Value = 1678.7778
$ws.Range("J132").Value = 3969.9333
$ws.Range("K132").Value = 5036.3334
$ws.Range("L132").Value = 11909.7999
$ws.Range("M132").Value = -2506.3334
$ws.Range("N132").Value = -16969.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 14000.272
$ws.Range("I132").Value = 4400.8
$ws.Range("J132").Value = 21999.834
$ws.Range("K132").Value = 13202.4
$ws.Range("L132").Value = 65999.50199999999
$ws.Range("M132").Value = -10672.4
$ws.Range("N132").Value = -71059.50199999999
